$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.518.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.425.56'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.624'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.86%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.423.46'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('E10').Value = '  +2.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.98'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.017.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.16'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.379.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000173'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.426.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '369.54'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.61'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('E25').Value = '  +6.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.534'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.68%  '
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.39'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.55'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.51'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.47%  '
$ws.Range('E38').Value = '  -1.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.81'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.63'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.43'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.03%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.721.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '335.97'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.92%  '
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '32.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.14%  '
$ws.Range('E51').Value = '  +3.91%  '

Write-Output "Applied cryptos update"
